$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = 'maa://24702 (94.61), maa://25390 (95.91), maa://36681 (87.34)'
$ws.Range('T2').Value = 'maa://22742 (90.7), *maa://20791 (62.82)'
$ws.Range('AB2').Value = 'maa://21246 (91.47), maa://36684 (95.83), ***maa://22731 (6.25)'
$ws.Range('D3').Value = 'maa://40192 (96.77), maa://36987 (96.15), maa://39849 (88.89)'
$ws.Range('T3').Value = 'maa://24617 (89.92), **maa://20790 (43.48), ***maa://37170 (16.18), maa://45854 (85.29)'
$ws.Range('X4').Value = '**maa://32495 (48.89), ***maa://31785 (22.22), maa://43217 (91.03), ***maa://36683 (28.26)'
$ws.Range('AE4').Value = '4'
$ws.Range('AF4').Value = '*maa://30062 (64.0), ***maa://26209 (13.04), *maa://39394 (65.38), maa://48095 (100.0)'
$ws.Range('D6').Value = 'maa://42407 (95.16)'
$ws.Range('L7').Value = 'maa://28624 (92.98), maa://24957 (97.78)'
$ws.Range('X7').Value = 'maa://22399 (95.68), *maa://22758 (75.34)'
$ws.Range('A8').Value = '更新日期：2025.04.03 13:19:19'
$ws.Range('L9').Value = 'maa://22762 (92.47), maa://39552 (81.25)'
$ws.Range('T10').Value = 'maa://27395 (96.59), maa://22755 (87.83), **maa://22756 (40.91), ***maa://21737 (10.61)'
$ws.Range('X10').Value = 'maa://22301 (97.78), maa://45828 (88.24), maa://22726 (100.0)'
$ws.Range('D13').Value = 'maa://24999 (92.18), maa://36673 (92.41), maa://25001 (85.92)'
$ws.Range('H13').Value = '*maa://21248 (73.9), **maa://22728 (46.67)'
$ws.Range('X13').Value = 'maa://34957 (81.93), **maa://22768 (50.0)'
$ws.Range('AF13').Value = '**maa://22737 (34.25), maa://39883 (90.79), *maa://39885 (51.61)'
$ws.Range('L14').Value = 'maa://26245 (96.82), maa://21288 (96.3), maa://39841 (94.17), maa://36682 (97.44)'
$ws.Range('T15').Value = 'maa://23892 (96.34)'
$ws.Range('AF15').Value = 'maa://21364 (80.94), *maa://36666 (77.12), *maa://22766 (68.33)'
$ws.Range('H17').Value = 'maa://22430 (88.83), maa://39599 (84.21)'
$ws.Range('H18').Value = 'maa://24421 (88.33)'
$ws.Range('L18').Value = 'maa://22466 (90.91), *maa://22732 (51.55)'
$ws.Range('T19').Value = 'maa://24386 (99.19)'
$ws.Range('AB19').Value = '*maa://30709 (66.22), *maa://36668 (57.5)'
$ws.Range('D20').Value = 'maa://21432 (90.61), maa://25198 (93.75), *maa://20795 (50.77), maa://36680 (91.18)'
$ws.Range('L20').Value = 'maa://41331 (85.03)'
$ws.Range('P21').Value = 'maa://24381 (80.95)'
$ws.Range('AF21').Value = 'maa://22524 (93.42), *maa://22432 (78.82)'
$ws.Range('X22').Value = 'maa://21282 (98.64), *maa://37649 (65.52)'
$ws.Range('D23').Value = '***maa://28036 (28.0), *maa://41753 (52.38)'
$ws.Range('L23').Value = 'maa://39756 (95.92), maa://39875 (94.59)'
$ws.Range('X24').Value = 'maa://29988 (83.71), maa://23504 (93.35), **maa://22892 (40.54), *maa://25141 (77.1), *maa://36663 (77.5), ***maa://22815 (23.08)'
$ws.Range('D25').Value = 'maa://29753 (95.26)'
$ws.Range('H25').Value = '*maa://29063 (72.78), *maa://25311 (74.77), ***maa://22725 (4.76), *maa://45047 (66.67)'
$ws.Range('X25').Value = '*maa://29890 (80.0)'
$ws.Range('AB26').Value = 'maa://42235 (94.78)'
$ws.Range('D28').Value = 'maa://24465 (90.98), maa://25725 (84.27)'
$ws.Range('X28').Value = 'maa://39929 (90.69), maa://41749 (91.67), ***maa://39723 (13.89)'
$ws.Range('H29').Value = '*maa://25175 (65.45)'
$ws.Range('AF29').Value = '*maa://24080 (68.93), maa://42865 (81.43), ***maa://34960 (8.33)'
$ws.Range('D30').Value = 'maa://45792 (94.44)'
$ws.Range('AB30').Value = 'maa://42979 (97.14), maa://45822 (100.0), *maa://45045 (80.0)'
$ws.Range('L31').Value = 'maa://35926 (93.4), maa://36258 (85.12), *maa://43904 (72.73)'
$ws.Range('H32').Value = 'maa://21895 (97.1), maa://36667 (97.73), **maa://20793 (38.78), maa://22760 (100.0)'
$ws.Range('L32').Value = 'maa://28065 (95.74)'
$ws.Range('P34').Value = 'maa://48817 (87.5)'
$ws.Range('AF34').Value = '*maa://32650 (76.19)'
$ws.Range('L35').Value = 'maa://41296 (96.51)'
$ws.Range('T35').Value = 'maa://24842 (94.23)'
$ws.Range('L37').Value = 'maa://45718 (97.66), *maa://47069 (75.0), maa://45789 (100.0)'
$ws.Range('AF38').Value = 'maa://36697 (86.94)'
$ws.Range('H39').Value = 'maa://36670 (89.22), maa://25199 (84.96), maa://30434 (91.3), *maa://45059 (78.26), ***maa://25036 (18.52), *maa://44165 (66.67)'
$ws.Range('T39').Value = '*maa://45788 (80.0), maa://47079 (93.33), *maa://45790 (73.33)'
$ws.Range('H46').Value = 'maa://35931 (91.98), maa://43901 (93.75)'
$ws.Range('H53').Value = 'maa://32534 (94.25), **maa://32434 (33.33)'
$ws.Range('H55').Value = 'maa://32532 (92.14)'
